$d = $word.ActiveDocument

# --- Locate the target paragraph (the one mentioning "OneR") ---------------
$count = $d.Paragraphs.Count
$target = $null
for ($i = 1; $i -le $count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*OneR*") {
        $target = $cand
        break
    }
}

if ($target -ne $null) {
    $rng = $target.Range
    $xml = $rng.WordOpenXML

    # 1. Drop the paragraph-mark run properties (<w:rPr>...</w:rPr>) that sit
    #    directly inside this paragraph's <w:pPr>, right before </w:pPr>.
    $xml = $xml -replace '(<w:pPr>.*?)<w:rPr>.*?</w:rPr></w:pPr>', '$1</w:pPr>'

    # 2. Unwrap the hyperlink around "OneR": keep its run (and wrap it with
    #    spell-check proofErr markers) but drop the <w:hyperlink> tags.
    $xml = $xml -replace '<w:hyperlink[^>]*>(<w:r[^>]*><w:rPr>.*?</w:rPr><w:t>OneR</w:t></w:r>)</w:hyperlink>', '<w:proofErr w:type="spellStart"/>$1<w:proofErr w:type="spellEnd"/>'

    $rng.InsertXML($xml)
}

# --- Remove the trailing empty paragraph at the end of the document body ---
$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)
if ($last.Range.Text -eq [string][char]13) {
    $prev = $d.Paragraphs.Item($count - 1)
    $d.Range($prev.Range.End, $last.Range.End).Delete()
}
